$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 1485938350
$ws.Range("B3").Value = 3
$ws.Range("A4").Value = 1485938350
$ws.Range("B4").Value = 4

$ws.Range("B5").Select()
